$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 668 (Excel shifts existing rows 668:725 down to 669:726,
# dimension grows from A1:T725 to A1:T726, formatting/styles of the row carry over).
$ws.Rows(668).Insert()

# Populate the newly inserted row with its data.
$newDate = Get-Date -Year 2023 -Month 8 -Day 28 -Hour 0 -Minute 0 -Second 0

$ws.Cells.Item(668, 1).Value = 6
$ws.Cells.Item(668, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(668, 3).Value = "Metropolitana"
$ws.Cells.Item(668, 4).Value = $newDate
$ws.Cells.Item(668, 5).Value = 13
$ws.Cells.Item(668, 6).Value = "Fruta"
$ws.Cells.Item(668, 7).Value = 100101
$ws.Cells.Item(668, 8).Value = "Berries"
$ws.Cells.Item(668, 9).Value = 100101001
$ws.Cells.Item(668, 10).Value = "Arándano (blue)"
$ws.Cells.Item(668, 11).Value = "Sin especificar"
$ws.Cells.Item(668, 12).Value = "Primera"
$ws.Cells.Item(668, 13).Value = 1070
$ws.Cells.Item(668, 14).Value = 12000
$ws.Cells.Item(668, 15).Value = 12000
$ws.Cells.Item(668, 16).Value = 12000
$ws.Cells.Item(668, 17).Value = "`$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(668, 18).Value = "Perú"
$ws.Cells.Item(668, 19).Value = 8000
$ws.Cells.Item(668, 20).Value = 1.5
